$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.983.67"
$ws.Range("E2").Value = "  +4.83%  "

$ws.Range("D3").Value = "4.041.66"
$ws.Range("E3").Value = "  +4.70%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.97"
$ws.Range("E5").Value = "  +2.78%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.64"
$ws.Range("E6").Value = "  +8.93%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.696"
$ws.Range("E7").Value = "  +14.26%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.765"
$ws.Range("E9").Value = "  +7.44%  "

$ws.Range("E10").Value = "  +4.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000333"
$ws.Range("E11").Value = "  +3.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.73"
$ws.Range("E12").Value = "  +16.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.95"
$ws.Range("E13").Value = "  +5.66%  "

$ws.Range("D14").Value = "4.684.07"
$ws.Range("E14").Value = "  +4.27%  "

$ws.Range("D15").Value = "4.053.58"
$ws.Range("E15").Value = "  +4.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.42"
$ws.Range("E16").Value = "  +2.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.73"
$ws.Range("E17").Value = "  -3.19%  "

$ws.Range("E18").Value = "  +1.59%  "

$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("D20").Value = "71.966.43"
$ws.Range("E20").Value = "  +4.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "436.51"
$ws.Range("E21").Value = "  +4.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "99.66"
$ws.Range("E22").Value = "  +14.67%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.59"
$ws.Range("E23").Value = "  +1.57%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.27"
$ws.Range("E24").Value = "  +6.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.72"
$ws.Range("E25").Value = "  +4.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.30"
$ws.Range("E26").Value = "  -4.64%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.00"
$ws.Range("E27").Value = "  +5.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.71"
$ws.Range("E28").Value = "  +30.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.18"
$ws.Range("E29").Value = "  +4.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.84"
$ws.Range("E30").Value = "  +2.73%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.75"
$ws.Range("E31").Value = "  +2.09%  "

$ws.Range("E32").Value = "  +6.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "682.35"
$ws.Range("E33").Value = "  +1.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.97"
$ws.Range("E34").Value = "  +1.00%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "67.34"
$ws.Range("E35").Value = "  +0.83%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "43.06"
$ws.Range("E36").Value = "  +9.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.436"
$ws.Range("E37").Value = "  -1.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.158"
$ws.Range("E38").Value = "  +6.59%  "

$ws.Range("D39").Value = "0.0₃0848"
$ws.Range("E39").Value = "  -0.41%  "

$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.51"
$ws.Range("E40").Value = "  +11.00%  "

$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.46"
$ws.Range("E41").Value = "  -1.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.11%  "

$ws.Range("E43").Value = "  +4.40%  "

$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("E45").Value = "  +8.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.73"
$ws.Range("E46").Value = "  -4.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.42"
$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.54"
$ws.Range("E48").Value = "  +8.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.07"
$ws.Range("E49").Value = "  +3.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.37"
$ws.Range("E50").Value = "  +2.78%  "

$ws.Range("E51").Value = "  -5.89%  "
